# Invoice template data update: vendor "MINDSHERPA" with three invoice numbers,
# replacing the old two placeholder rows and adding a new third data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep the invoice-number column formatted as text so values like "85420"
# are not coerced into numbers (and stay as shared-string entries).
$ws.Range("A1:B4").NumberFormat = "@"

$ws.Range("A2").Value = "MINDSHERPA"
$ws.Range("B2").Value = "85420"

$ws.Range("A3").Value = "MINDSHERPA"
$ws.Range("B3").Value = "851244"

$ws.Range("A4").Value = "MINDSHERPA"
$ws.Range("B4").Value = "850888"

# Tidy up the page margins to the workbook defaults.
$ws.PageSetup.LeftMargin = $excel.InchesToPoints(0.75)
$ws.PageSetup.RightMargin = $excel.InchesToPoints(0.75)
$ws.PageSetup.TopMargin = $excel.InchesToPoints(0.75)
$ws.PageSetup.BottomMargin = $excel.InchesToPoints(0.5)
$ws.PageSetup.HeaderMargin = $excel.InchesToPoints(0.5)
$ws.PageSetup.FooterMargin = $excel.InchesToPoints(0.75)

# Collapse the selection down to the single anchor cell.
$null = $ws.Range("A1").Select()
